# Fix branch-code trimming issue: Branch Code values in column F were
# stored without their leading "1" (or "11" for codes starting with "2"),
# e.g. "0219110" should really be "100219110" and "2124610" should
# really be "112124610". Re-pad every Branch Code in column F so the
# codes are always 9 digits long, per the updated pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row   # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 355 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)   # column F = Branch Code
    $old = [string]$cell.Value2

    if ([string]::IsNullOrEmpty($old)) { continue }

    if ($old.Length -eq 7 -and ($old.StartsWith("0") -or $old.StartsWith("1") -or $old.StartsWith("2"))) {
        if ($old.StartsWith("2")) {
            $prefix = "11"
        } else {
            $prefix = "10"
        }
        # NOTE: use string interpolation (not "+") to concatenate, since
        # "+" on numeric-looking strings performs arithmetic addition.
        $new = "$prefix$old"

        # The new code is fully numeric-looking ("100219110"), so a plain
        # Value2 assignment would store it as a number and lose the
        # original text semantics (F column is plain text in the source
        # file). Mark the cell as text first so the string is kept as
        # text, then clear the leftover "@" number-format style so no
        # stray formatting is introduced.
        $cell.NumberFormat = "@"
        $cell.Value2 = $new
        $cell.ClearFormats()
    }
}
